# Applies the weekly refresh described in the commit:
# "Fruta / hortaliza, semanal" — existing price rows for
# Agricola del Norte S.A. de Arica / Pera shift to the next
# reporting week, and two new rows (22-23) are appended.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the cells that changed value in rows 10-21 ---
# Row 10
$ws.Range("D10").Value = 44497
$ws.Range("M10").Value = 300
$ws.Range("Q10").Value = '$/bandeja 18 kilos granel'
$ws.Range("R10").Value = 'Región de O''Higgins'
# Row 11
$ws.Range("Q11").Value = '$/bandeja 18 kilos granel'
$ws.Range("R11").Value = 'Región de O''Higgins'
# Row 12
$ws.Range("D12").Value = 44355
$ws.Range("K12").Value = 'Packham''s Triumph'
$ws.Range("M12").Value = 200
$ws.Range("N12").Value = 17000
$ws.Range("O12").Value = 18000
$ws.Range("P12").Value = 17500
$ws.Range("Q12").Value = '$/caja 18 kilos granel'
$ws.Range("R12").Value = 'Región Metropolitana'
$ws.Range("S12").Value = 972
# Row 13
$ws.Range("D13").Value = 44355
$ws.Range("K13").Value = 'Winter Nelis'
$ws.Range("M13").Value = 250
$ws.Range("N13").Value = 17000
$ws.Range("O13").Value = 18000
$ws.Range("P13").Value = 17500
$ws.Range("Q13").Value = '$/caja 18 kilos granel'
$ws.Range("R13").Value = 'Región Metropolitana'
$ws.Range("S13").Value = 972
# Row 14
$ws.Range("D14").Value = 44336
$ws.Range("N14").Value = 21000
$ws.Range("O14").Value = 22000
$ws.Range("P14").Value = 21500
$ws.Range("S14").Value = 1194
# Row 15
$ws.Range("D15").Value = 44421
$ws.Range("M15").Value = 270
$ws.Range("N15").Value = 16000
$ws.Range("O15").Value = 17000
$ws.Range("P15").Value = 16500
$ws.Range("Q15").Value = '$/bandeja 18 kilos granel'
$ws.Range("S15").Value = 917
# Row 16
$ws.Range("D16").Value = 44421
$ws.Range("N16").Value = 16000
$ws.Range("O16").Value = 17000
$ws.Range("P16").Value = 16500
$ws.Range("Q16").Value = '$/bandeja 18 kilos granel'
$ws.Range("S16").Value = 917
# Row 17
$ws.Range("D17").Value = 44292
$ws.Range("M17").Value = 300
$ws.Range("N17").Value = 22000
$ws.Range("O17").Value = 23000
$ws.Range("P17").Value = 22500
$ws.Range("S17").Value = 1250
# Row 18
$ws.Range("D18").Value = 44292
$ws.Range("M18").Value = 250
$ws.Range("N18").Value = 22000
$ws.Range("O18").Value = 23000
$ws.Range("P18").Value = 22500
$ws.Range("S18").Value = 1250
# Row 19
$ws.Range("D19").Value = 44280
$ws.Range("M19").Value = 350
$ws.Range("N19").Value = 24000
$ws.Range("O19").Value = 25000
$ws.Range("P19").Value = 24500
$ws.Range("Q19").Value = '$/caja 18 kilos granel'
$ws.Range("S19").Value = 1361
# Row 20
$ws.Range("D20").Value = 44280
$ws.Range("M20").Value = 300
$ws.Range("N20").Value = 24000
$ws.Range("O20").Value = 25000
$ws.Range("P20").Value = 24500
$ws.Range("Q20").Value = '$/caja 18 kilos granel'
$ws.Range("S20").Value = 1361
# Row 21
$ws.Range("D21").Value = 44474
$ws.Range("M21").Value = 270
$ws.Range("N21").Value = 18000
$ws.Range("O21").Value = 19000
$ws.Range("P21").Value = 18500
$ws.Range("Q21").Value = '$/caja 18 kilos empedrada'
$ws.Range("S21").Value = 1028

# --- Append the two new rows (22 and 23) ---
# Row 22
$ws.Range("A22").Value = 1
$ws.Range("B22").Value = 'Agrícola del Norte S.A. de Arica'
$ws.Range("C22").Value = 'Arica y Parinacota'
$ws.Range("D22").Value = 44474
$ws.Range("E22").Value = 15
$ws.Range("F22").Value = 'Fruta'
$ws.Range("G22").Value = 100104
$ws.Range("H22").Value = 'Frutos de pepita'
$ws.Range("I22").Value = 100104005
$ws.Range("J22").Value = 'Pera'
$ws.Range("K22").Value = 'Winter Nelis'
$ws.Range("L22").Value = 'Segunda'
$ws.Range("M22").Value = 250
$ws.Range("N22").Value = 17000
$ws.Range("O22").Value = 18000
$ws.Range("P22").Value = 17500
$ws.Range("Q22").Value = '$/bandeja 18 kilos granel'
$ws.Range("R22").Value = 'Región de O''Higgins'
$ws.Range("S22").Value = 972
$ws.Range("T22").Value = 18
$ws.Range("D22").NumberFormat = "YYYY-MM-DD HH:MM:SS"
# Row 23
$ws.Range("A23").Value = 1
$ws.Range("B23").Value = 'Agrícola del Norte S.A. de Arica'
$ws.Range("C23").Value = 'Arica y Parinacota'
$ws.Range("D23").Value = 44314
$ws.Range("E23").Value = 15
$ws.Range("F23").Value = 'Fruta'
$ws.Range("G23").Value = 100104
$ws.Range("H23").Value = 'Frutos de pepita'
$ws.Range("I23").Value = 100104005
$ws.Range("J23").Value = 'Pera'
$ws.Range("K23").Value = 'Packham''s Triumph'
$ws.Range("L23").Value = 'Segunda'
$ws.Range("M23").Value = 250
$ws.Range("N23").Value = 17000
$ws.Range("O23").Value = 18000
$ws.Range("P23").Value = 17500
$ws.Range("Q23").Value = '$/caja 18 kilos granel'
$ws.Range("R23").Value = 'Región de O''Higgins'
$ws.Range("S23").Value = 972
$ws.Range("T23").Value = 18
$ws.Range("D23").NumberFormat = "YYYY-MM-DD HH:MM:SS"
